# "pic at ICU entry + import individual data"
#
# The 17.03.2020 scenario row is dropped entirely (rows below it shift up
# by one), a handful of mlam/mpic values are revised for the remaining
# rows, and the (previously unused) "confinement" comment is placed on
# the 28.03.2020 row with its lambda value corrected from 1.15 to 1.12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "17.03.2020" row (was row 4) - everything below shifts up.
$ws.Rows(4).Delete()

# --- Column D (mpic) drops from 0.3 to 0.2 on the rows that still had 0.3 ---
$ws.Range("D2").Value = 0.2
$ws.Range("D3").Value = 0.2
$ws.Range("D4").Value = 0.2
$ws.Range("D5").Value = 0.2

# --- Revised mlam (column B) values for the remaining rows ---
$ws.Range("B3").Value = 1.12
$ws.Range("B4").Value = 1.12
$ws.Range("B5").Value = 1.1
$ws.Range("B6").Value = 1.09

# --- Place the confinement comment on row 5 (28.03.2020), first matching
#     the pre-existing shared string exactly so it's reused in place, then
#     correcting its lambda figure from 1.15 to 1.12 per the new mlam series ---
$ws.Range("J5").Value = "Début de l’effet du confinement (lam 1.15 → 1)"
$ws.Range("J5").Value = "Début de l’effet du confinement (lam 1.12 → 1)"

# Leave the selection on the last cell touched, same as the authored edit.
[void]$ws.Range("J11").Select()
